$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'NO FOS',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'NO FOS',
    'ANIL KUMAR PANDIT',
    'NO FOS',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'BEENESH',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'MANOJ TOMAR',
    'PRAMOD-JAIPUR',
    'MANOJ TOMAR',
    'NO FOS',
    'NO FOS',
    'SHAKEEL',
    'NO FOS',
    'MANOJ TOMAR',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'PRAMOD-JAIPUR',
    'ANIL KUMAR PANDIT',
    'MANOJ TOMAR',
    'NO FOS',
    'MANOJ TOMAR',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'SURESH',
    'MANOJ TOMAR',
    'MANOJ TOMAR',
    'NO FOS',
    'ANIL KUMAR PANDIT',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'MANOJ TOMAR',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'MANOJ TOMAR',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'RAJENDER KUMAR VERMA',
    'GIRJESH KUMAR',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'PRAMOD-JAIPUR',
    'SUKHVINDER SINGH',
    'NO FOS',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'SHAKEEL',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'NO FOS',
    'MANOJ TOMAR',
    'PRAMOD-JAIPUR',
    'PRAMOD-JAIPUR',
    'ANIL KUMAR PANDIT',
    'SUBHASH CHAND',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'SUKHVINDER SINGH',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'MANOJ TOMAR',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'ANIL KUMAR PANDIT',
    'NO FOS',
    'BEENESH',
    'NO FOS',
    'ANIL KUMAR PANDIT',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'SURESH',
    'SUKHVINDER SINGH',
    'MANOJ TOMAR',
    'BEENESH',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'BEENESH',
    'NO FOS',
    'BEENESH',
    'NO FOS',
    'MANOJ TOMAR',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'SUBHASH CHAND',
    'NO FOS',
    'SUKHVINDER SINGH',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'ANIL KUMAR PANDIT',
    'NO FOS',
    'SUBHASH CHAND',
    'BEENESH',
    'NO FOS',
    'NO FOS',
    'MANOJ TOMAR',
    'NO FOS',
    'NO FOS',
    'MANOJ TOMAR',
    'NO FOS',
    'NO FOS',
    'MANOJ TOMAR',
    'NO FOS',
    'SUBHASH CHAND',
    'MANOJ TOMAR',
    'NO FOS',
    'NO FOS',
    'MANOJ TOMAR',
    'NO FOS',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'ARVIND KUMAR',
    'BEENESH',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'SHAKEEL',
    'GURPREET SINGH',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'SHAKEEL',
    'SHAKEEL',
    'ANIL KUMAR PANDIT',
    'MANOJ TOMAR',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'RAJENDER KUMAR VERMA',
    'MANOJ TOMAR',
    'NO FOS',
    'NO FOS',
    'ANIL KUMAR PANDIT',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'SUKHVINDER SINGH',
    'NO FOS',
    'NO FOS',
    'SUBHASH CHAND',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'ARVIND KUMAR',
    'BEENESH',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'BEENESH',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'BEENESH',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'SURESH',
    'NO FOS',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'GURPREET SINGH',
    'NO FOS',
    'NO FOS',
    'SURESH',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'MANOJ TOMAR',
    'NO FOS',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'SUBHASH CHAND',
    'NO FOS',
    'NO FOS',
    'BEENESH',
    'BEENESH',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'BEENESH',
    'SUKHVINDER SINGH',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'SUKHVINDER SINGH',
    'ARVIND KUMAR',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'GURPREET SINGH',
    'NO FOS',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'ANIL KUMAR PANDIT',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'PRAMOD-JAIPUR',
    'MANOJ TOMAR',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'SUBHASH CHAND',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'BEENESH',
    'NO FOS',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'GIRJESH KUMAR',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'BEENESH',
    'NO FOS',
    'SUBHASH CHAND',
    'NO FOS',
    'ANIL KUMAR PANDIT',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'ANIL KUMAR PANDIT',
    'NO FOS',
    'SURESH',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'BEENESH',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'RAJENDER KUMAR VERMA',
    'NO FOS',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'ANIL KUMAR PANDIT',
    'NO FOS',
    'NO FOS',
    'RAJENDER KUMAR VERMA',
    'NO FOS',
    'SURESH',
    'NO FOS',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'BEENESH',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'SHAKEEL',
    'NO FOS',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'ANIL KUMAR PANDIT',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'SHAKEEL',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'PRAMOD-JAIPUR',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'NO FOS',
    'NO FOS'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $values[$i]
}
